# Apply updated "想去人数" (want-to-go count) values as published at
# gh-pages output generated at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 6026
$ws1.Range("F6").Value = 65
$ws1.Range("F19").Value = 4705
$ws1.Range("F20").Value = 117
$ws1.Range("F21").Value = 53
$ws1.Range("F29").Value = 51

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 78
$ws2.Range("F17").Value = 75
$ws2.Range("F23").Value = 507

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 6026
$ws4.Range("F28").Value = 4705
$ws4.Range("F29").Value = 53
$ws4.Range("F37").Value = 51
$ws4.Range("F43").Value = 507
